{"js": "// Word JS API (Office.js) script.\n// Goal (per diff):\n//  1. Convert the hyperlinked repo-URL run into a plain (non-hyperlinked) run,\n//     keeping the same visible text.\n//  2. Remove the \"Viikko 7\" paragraph, the blank paragraph that follows it,\n//     and the Finnish commentary paragraph (\"En saanut ollenkaan koodia ...\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate paragraphs by their text content so the script is resilient to\n// exact indices.\nlet hyperlinkParaIndex = -1;\nlet viikkoParaIndex = -1;\nlet finnishParaIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t.indexOf(\"github.com/LeoSuzu/Data_Structure_and_Algorythms\") !== -1) {\n    hyperlinkParaIndex = i;\n  } else if (t === \"Viikko 7\") {\n    viikkoParaIndex = i;\n  } else if (t.indexOf(\"En saanut ollenkaan koodia\") !== -1) {\n    finnishParaIndex = i;\n  }\n}\n\n// 1) Strip the hyperlink from the repo-URL paragraph, leaving plain text.\nif (hyperlinkParaIndex !== -1) {\n  const linkRange = items[hyperlinkParaIndex].getRange();\n  linkRange.hyperlink = \"\";\n}\n\n// 2) Delete the \"Viikko 7\" paragraph, the blank paragraph right after it,\n//    and the Finnish paragraph after that.\nif (viikkoParaIndex !== -1) {\n  const blankAfterViikkoIndex = viikkoParaIndex + 1;\n  items[viikkoParaIndex].delete();\n  if (blankAfterViikkoIndex < items.length) {\n    items[blankAfterViikkoIndex].delete();\n  }\n}\nif (finnishParaIndex !== -1) {\n  items[finnishParaIndex].delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Goal (per diff):\n#  1. Convert the hyperlinked repo-URL run into a plain (non-hyperlinked) run,\n#     keeping the same visible text and its en-US language tag.\n#  2. Remove the \"Viikko 7\" paragraph, the blank paragraph that follows it,\n#     and the Finnish commentary paragraph (\"En saanut ollenkaan koodia ...\").\n\n$d = $word.ActiveDocument\n\n# --- Locate the paragraphs we care about by their text, so this is not\n#     dependent on brittle fixed indices. ---\n$hyperlinkParaIndex = -1\n$viikkoParaIndex = -1\n$finnishParaIndex = -1\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -like \"*github.com/LeoSuzu/Data_Structure_and_Algorythms*\") {\n        $hyperlinkParaIndex = $p.Index\n    } elseif ($t -eq \"Viikko 7\") {\n        $viikkoParaIndex = $p.Index\n    } elseif ($t -like \"*En saanut ollenkaan koodia*\") {\n        $finnishParaIndex = $p.Index\n    }\n}\n\n# 1) Strip the hyperlink from the repo-URL paragraph, leaving plain text\n#    with the same language formatting it already had.\nif ($hyperlinkParaIndex -ne -1) {\n    $h = $d.Hyperlinks(1)\n    $linkText = $h.TextToDisplay\n    $origLanguageID = $h.Range.LanguageID\n    $h.Delete()\n\n    $linkPara = $d.Paragraphs($hyperlinkParaIndex).Range\n    $textOnly = $d.Range($linkPara.Start, $linkPara.End - 1)\n    $textOnly.Delete()\n\n    $ins = $d.Paragraphs($hyperlinkParaIndex).Range\n    $ins.InsertBefore($linkText)\n    $d.Paragraphs($hyperlinkParaIndex).Range.LanguageID = $origLanguageID\n}\n\n# 2) Delete the \"Viikko 7\" paragraph, the blank paragraph right after it,\n#    and the Finnish paragraph after that (three consecutive paragraphs).\nif ($viikkoParaIndex -ne -1) {\n    $d.Paragraphs($viikkoParaIndex).Range.Delete()\n    $d.Paragraphs($viikkoParaIndex).Range.Delete()\n    $d.Paragraphs($viikkoParaIndex).Range.Delete()\n}\n"}
